# Bump the "Förändrad" (Changed) date column (C) by one day for every
# data row (2-161): serial date 46074 (2026-02-21) -> 46075 (2026-02-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 161 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46074) {
        $cell.Value2 = 46075
    }
}
